$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split every "left-bottom" run into "left-" (keeps the original run
#    formatting) followed by a new "middle" run that carries only the
#    <w:rtl w:val="0"/> property (no explicit color), matching the target
#    diff. We locate a donor range elsewhere in the document whose run
#    formatting already consists of nothing but <w:rtl w:val="0"/> and use
#    Range.FormattedText to transplant that exact (minimal) run formatting
#    onto the new text, then rename the transplanted text to "middle".
# ---------------------------------------------------------------------------

function Split-LeftBottom($searchStart) {
    $tailOfDoc = $d.Range($searchStart, $d.Content.End)
    $found = $tailOfDoc.Find.Execute("left-bottom", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }

    $matchStart = $tailOfDoc.Start
    $matchEnd = $tailOfDoc.End

    # Donor run: the standalone "mie" inside "<m>mie</m>" has rPr = { rtl=0 }
    # only (no color), which is exactly the formatting the new "middle" run
    # needs.
    $donor = $d.Content
    $donor.Find.Execute("mie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $donorLen = $donor.End - $donor.Start

    # Replace the "bottom" tail (keeping "left-" untouched) with a copy of
    # the donor's run formatting.
    $tail = $d.Range($matchStart + 5, $matchEnd)
    $tail.FormattedText = $donor.FormattedText

    # Rename the freshly-formatted text from "mie" to "middle" in place.
    $newRun = $d.Range($matchStart + 5, $matchStart + 5 + $donorLen)
    $newRun.Text = "middle"

    return $matchStart + 5 + 6
}

$pos = 0
while ($true) {
    $pos = Split-LeftBottom $pos
    if ($pos -lt 0) { break }
}

# ---------------------------------------------------------------------------
# 2) Give the section a footer margin of 720 twips (36 pt), i.e. add
#    w:footer="720" to <w:pgMar>.
# ---------------------------------------------------------------------------
$d.PageSetup.FooterDistance = 36
